$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B/C/D/E in this sheet are stored as text (inline strings) in the
# source data, including numeric-looking price figures. When a plain numeric
# string is assigned to a cell, Excel auto-converts it to a Number, so any
# target value that parses as a number needs the cell pre-formatted as Text.

$ws.Range("D2").Value = '27.094.95'
$ws.Range("E2").Value = '  -2.60%  '

$ws.Range("D3").Value = '1.732.59'
$ws.Range("E3").Value = '  -1.39%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9992'
$ws.Range("E4").Value = '  -0.32%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.66'
$ws.Range("E5").Value = '  -5.27%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9992'
$ws.Range("E6").Value = '  -0.21%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4859'
$ws.Range("E7").Value = '  +6.15%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3513'
$ws.Range("E8").Value = '  +0.66%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '42.14'
$ws.Range("E9").Value = '  +0.18%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07292'
$ws.Range("E10").Value = '  -0.64%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.056'
$ws.Range("E11").Value = '  -2.54%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9989'
$ws.Range("E12").Value = '  -0.29%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.05'
$ws.Range("E13").Value = '  -2.53%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.908'
$ws.Range("E14").Value = '  -0.98%  '

$ws.Range("D15").Value = '1.726.77'
$ws.Range("E15").Value = '  -1.73%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.900'
$ws.Range("E16").Value = '  -3.62%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '87.37'
$ws.Range("E17").Value = '  -4.68%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001041'
$ws.Range("E18").Value = '  -0.79%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06410'
$ws.Range("E19").Value = '  -0.17%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9991'
$ws.Range("E20").Value = '  -0.25%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.58'
$ws.Range("E21").Value = '  -1.54%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.690'
$ws.Range("E22").Value = '  -0.62%  '

$ws.Range("D23").Value = '27.140.66'
$ws.Range("E23").Value = '  -2.56%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.87'
$ws.Range("E24").Value = '  -2.51%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.082'
$ws.Range("E25").Value = '  -3.31%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '153.39'
$ws.Range("E26").Value = '  -5.36%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.02'
$ws.Range("E27").Value = '  +0.21%  '

$ws.Range("D28").Value = '1.926.19'
$ws.Range("E28").Value = '  -1.73%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.092'
$ws.Range("E29").Value = '  -2.75%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '122.06'
$ws.Range("E30").Value = '  -0.81%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.044'
$ws.Range("E31").Value = '  -3.46%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09333'
$ws.Range("E32").Value = '  +0.40%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.606'
$ws.Range("E33").Value = '  -1.01%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.421'
$ws.Range("E34").Value = '  -1.87%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02206'
$ws.Range("E35").Value = '  -1.98%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.05939'
$ws.Range("E36").Value = '  -2.35%  '

# Row 37: Aptos
$ws.Range("B37").Value = 'Aptos'
$ws.Range("C37").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '11.05'
$ws.Range("E37").Value = '  -5.72%  '

# Row 38: WEMIXTOKEN
$ws.Range("B38").Value = 'WEMIXTOKEN'
$ws.Range("C38").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.435'
$ws.Range("E38").Value = '  +4.94%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2008'
$ws.Range("E39").Value = '  -2.76%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.787'
$ws.Range("E40").Value = '  -2.11%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6030'
$ws.Range("E41").Value = '  -2.49%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9986'
$ws.Range("E42").Value = '  -0.03%  '

$ws.Range("E43").Value = '  -7.09%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.531'
$ws.Range("E44").Value = '  -2.54%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '12.81'

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.586'
$ws.Range("E46").Value = '  -3.61%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5682'
$ws.Range("E47").Value = '  -1.73%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '118.84'
$ws.Range("E48").Value = '  -2.33%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.851'
$ws.Range("E49").Value = '  -3.59%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.111'
$ws.Range("E50").Value = '  -0.88%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06655'
$ws.Range("E51").Value = '  -1.81%  '
